$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1) to snake_case field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize "de"/"y" in selected municipality names
$ws.Range("B2").Value = "Comitán De Domínguez"
$ws.Range("B6").Value = "Mazapa De Madero"
$ws.Range("B18").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B22").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B25").Value = "Villa De Etla"
$ws.Range("B30").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B44").Value = "Cosamaloapan De Carpio"
$ws.Range("B45").Value = "Martínez De La Torre"

# Remove trailing footer/metadata rows (rows 51-55)
$ws.Range("A51:A55").EntireRow.Delete()
